# Scheduled-runner refresh: update computed profit columns (H:N) across Sheets
# per the latest Kujata pricing pull. Values only -- no formulas/formatting involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 45
$ws.Range("H45").Value = 187
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = $null

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 10029.546
$ws.Range("I2").Value = 912.375
$ws.Range("K2").Value = 912.375
$ws.Range("M2").Value = -799.375

# Row 63
$ws.Range("H63").Value = 2362.5
$ws.Range("I63").Value = 2357.1428
$ws.Range("J63").Value = 2400
$ws.Range("K63").Value = 2357.1428
$ws.Range("L63").Value = 2400
$ws.Range("M63").Value = -1671.1428
$ws.Range("N63").Value = -3772

# Row 66
$ws.Range("H66").Value = 2362.5
$ws.Range("I66").Value = 2357.1428
$ws.Range("J66").Value = 2400
$ws.Range("K66").Value = 11785.714
$ws.Range("L66").Value = 12000
$ws.Range("M66").Value = -8353.714
$ws.Range("N66").Value = -18864

# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = $null
$ws.Range("N68").Value = 0

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = $null
$ws.Range("N71").Value = 0

# Row 74
$ws.Range("H74").Value = 1166.375
$ws.Range("I74").Value = 838.5454999999999
$ws.Range("J74").Value = 1887.6
$ws.Range("K74").Value = 838.5454999999999
$ws.Range("L74").Value = 1887.6
$ws.Range("M74").Value = 35.45450000000005
$ws.Range("N74").Value = -3635.6

# Row 77
$ws.Range("H77").Value = 1166.375
$ws.Range("I77").Value = 838.5454999999999
$ws.Range("J77").Value = 1887.6
$ws.Range("K77").Value = 4192.7275
$ws.Range("L77").Value = 9438
$ws.Range("M77").Value = 175.2725
$ws.Range("N77").Value = -18174

# Row 97
$ws.Range("H97").Value = 580
$ws.Range("I97").Value = 396
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 396
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = 100
$ws.Range("N97").Value = -2492

# Row 116
$ws.Range("H116").Value = 10029.546
$ws.Range("I116").Value = 912.375
$ws.Range("K116").Value = 912.375
$ws.Range("M116").Value = 1381.625

# Row 133
$ws.Range("H133").Value = 28619.715
$ws.Range("I133").Value = 28000
$ws.Range("J133").Value = 28723
$ws.Range("K133").Value = 28000
$ws.Range("L133").Value = 28723
$ws.Range("M133").Value = -25470
$ws.Range("N133").Value = -33783

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 10029.546
$ws.Range("I3").Value = 912.375
$ws.Range("K3").Value = 912.375
$ws.Range("M3").Value = -798.375

# Row 94
$ws.Range("H94").Value = 19231628
$ws.Range("I94").Value = 20834096
$ws.Range("K94").Value = 20834096
$ws.Range("M94").Value = -20833645

# Row 107
$ws.Range("H107").Value = 1874.1578
$ws.Range("I107").Value = 1353.3
$ws.Range("J107").Value = 2452.889
$ws.Range("K107").Value = 1353.3
$ws.Range("L107").Value = 2452.889
$ws.Range("M107").Value = 566.7
$ws.Range("N107").Value = -6292.889

$ws = $wb.Worksheets.Item("CRP")
# Row 93
$ws.Range("H93").Value = 17223
$ws.Range("J93").Value = 52500
$ws.Range("L93").Value = 52500
$ws.Range("N93").Value = -56244

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 493.2
$ws.Range("I7").Value = 516
$ws.Range("J7").Value = 402
$ws.Range("K7").Value = 1548
$ws.Range("L7").Value = 1206
$ws.Range("M7").Value = -1436
$ws.Range("N7").Value = -1430

# Row 34
$ws.Range("H34").Value = 1581.0526
$ws.Range("J34").Value = 2000
$ws.Range("L34").Value = 6000
$ws.Range("N34").Value = -6168

# Row 39
$ws.Range("H39").Value = 1898.2222
$ws.Range("J39").Value = 1898.2222
$ws.Range("L39").Value = 5694.6666
$ws.Range("N39").Value = -6282.6666

# Row 55
$ws.Range("H55").Value = 3047.2727
$ws.Range("J55").Value = 3047.2727
$ws.Range("L55").Value = 9141.8181
$ws.Range("N55").Value = -9495.8181

# Row 63
$ws.Range("H63").Value = 4008.5386
$ws.Range("I63").Value = 885.1667
$ws.Range("K63").Value = 2655.5001
$ws.Range("M63").Value = -1906.5001

# Row 66
$ws.Range("H66").Value = 4008.5386
$ws.Range("I66").Value = 885.1667
$ws.Range("K66").Value = 7966.5003
$ws.Range("M66").Value = -4222.5003

# Row 69
$ws.Range("H69").Value = 2286.6428
$ws.Range("I69").Value = 1299.5
$ws.Range("J69").Value = 2451.1667
$ws.Range("K69").Value = 3898.5
$ws.Range("L69").Value = 7353.500100000001
$ws.Range("M69").Value = -3087.5
$ws.Range("N69").Value = -8975.500100000001

# Row 72
$ws.Range("H72").Value = 2286.6428
$ws.Range("I72").Value = 1299.5
$ws.Range("J72").Value = 2451.1667
$ws.Range("K72").Value = 11695.5
$ws.Range("L72").Value = 22060.5003
$ws.Range("M72").Value = -7639.5
$ws.Range("N72").Value = -30172.5003

# Row 86
$ws.Range("H86").Value = 1033
$ws.Range("I86").Value = 1033
$ws.Range("K86").Value = 3099
$ws.Range("M86").Value = -1913

# Row 87
$ws.Range("H87").Value = 2778.5
$ws.Range("I87").Value = 814
$ws.Range("J87").Value = 3433.3333
$ws.Range("K87").Value = 2442
$ws.Range("L87").Value = 10299.9999
$ws.Range("M87").Value = -1194
$ws.Range("N87").Value = -12795.9999

# Row 89
$ws.Range("H89").Value = 1033
$ws.Range("I89").Value = 1033
$ws.Range("K89").Value = 9297
$ws.Range("M89").Value = -3369

# Row 90
$ws.Range("H90").Value = 2778.5
$ws.Range("I90").Value = 814
$ws.Range("J90").Value = 3433.3333
$ws.Range("K90").Value = 7326
$ws.Range("L90").Value = 30899.9997
$ws.Range("M90").Value = -1086
$ws.Range("N90").Value = -43379.9997

# Row 107
$ws.Range("H107").Value = 1100.5714
$ws.Range("I107").Value = 403
$ws.Range("J107").Value = 1216.8334
$ws.Range("K107").Value = 1209
$ws.Range("L107").Value = 3650.5002
$ws.Range("M107").Value = 711
$ws.Range("N107").Value = -7490.5002

# Row 122
$ws.Range("H122").Value = 885.375
$ws.Range("I122").Value = 826.6
$ws.Range("J122").Value = 983.3333
$ws.Range("K122").Value = 7439.400000000001
$ws.Range("L122").Value = 8849.9997
$ws.Range("M122").Value = -4989.400000000001
$ws.Range("N122").Value = -13749.9997

# Row 131
$ws.Range("H131").Value = 37038708
$ws.Range("I131").Value = 142857420
$ws.Range("J131").Value = 2153.95
$ws.Range("K131").Value = 428572260
$ws.Range("L131").Value = 6461.849999999999
$ws.Range("M131").Value = -428567220
$ws.Range("N131").Value = -16541.85

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 21434748
$ws.Range("I70").Value = 19236094
$ws.Range("J70").Value = 25007562
$ws.Range("K70").Value = 19236094
$ws.Range("L70").Value = 25007562
$ws.Range("M70").Value = -19235824
$ws.Range("N70").Value = -25008102

# Row 73
$ws.Range("H73").Value = 21434748
$ws.Range("I73").Value = 19236094
$ws.Range("J73").Value = 25007562
$ws.Range("K73").Value = 19236094
$ws.Range("L73").Value = 25007562
$ws.Range("M73").Value = -19235158
$ws.Range("N73").Value = -25009434

# Row 104
$ws.Range("H104").Value = 38534.2
$ws.Range("J104").Value = 38534.2
$ws.Range("L104").Value = 38534.2
$ws.Range("N104").Value = -45522.2

# Row 113
$ws.Range("H113").Value = 1979.68
$ws.Range("I113").Value = 1209.2727
$ws.Range("J113").Value = 2585
$ws.Range("K113").Value = 1209.2727
$ws.Range("L113").Value = 2585
$ws.Range("M113").Value = 960.7273
$ws.Range("N113").Value = -6925

# Row 132
$ws.Range("H132").Value = 3066.5386
$ws.Range("I132").Value = 2486.7
$ws.Range("K132").Value = 7460.099999999999
$ws.Range("M132").Value = -4930.099999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1147.9333
$ws.Range("I16").Value = 1070.091
$ws.Range("J16").Value = 1362
$ws.Range("K16").Value = 1070.091
$ws.Range("L16").Value = 1362
$ws.Range("M16").Value = -900.0909999999999
$ws.Range("N16").Value = -1702

# Row 93
$ws.Range("H93").Value = 692.25
$ws.Range("I93").Value = 692.25
$ws.Range("K93").Value = 692.25
$ws.Range("M93").Value = 555.75

# Row 132
$ws.Range("H132").Value = 79831.46000000001
$ws.Range("J132").Value = 127750.875
$ws.Range("L132").Value = 383252.625
$ws.Range("N132").Value = -388312.625

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3014.3076
$ws.Range("I132").Value = 2426.1365
$ws.Range("J132").Value = 6249.25
$ws.Range("K132").Value = 7278.4095
$ws.Range("L132").Value = 18747.75
$ws.Range("M132").Value = -4748.4095
$ws.Range("N132").Value = -23807.75

# Row 136
$ws.Range("H136").Value = 652.82355
$ws.Range("I136").Value = 428.3
$ws.Range("K136").Value = 1284.9
$ws.Range("M136").Value = 1265.1
